$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.173.21'
$ws.Range("E2").Value = '  +2.97%  '

$ws.Range("D3").Value = '2.653.16'
$ws.Range("E3").Value = '  +2.71%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = "'595.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.28%  '

$ws.Range("D6").Value = "'156.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.82%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = "'0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.26%  '

$ws.Range("D9").Value = "'0.119"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.83%  '

$ws.Range("E10").Value = '  +4.41%  '

$ws.Range("E11").Value = '  +1.68%  '

$ws.Range("E12").Value = '  +1.86%  '

$ws.Range("D13").Value = "'29.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +5.55%  '

$ws.Range("D14").Value = "'0.0000189"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +21.62%  '

$ws.Range("D15").Value = '3.128.36'
$ws.Range("E15").Value = '  +2.67%  '

$ws.Range("D16").Value = '65.059.45'
$ws.Range("E16").Value = '  +3.09%  '

$ws.Range("D17").Value = '2.590.85'
$ws.Range("E17").Value = '  -1.05%  '

$ws.Range("D18").Value = "'12.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.36%  '

$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.84%  '

$ws.Range("D20").Value = "'355.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.96%  '

$ws.Range("D21").Value = "'7.30"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.53%  '

$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").Value = "'68.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("D25").Value = "'9.54"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.30%  '

$ws.Range("E26").Value = '  -0.98%  '

$ws.Range("D27").Value = "'8.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.44%  '

$ws.Range("E28").Value = '  +2.12%  '

$ws.Range("D29").Value = '0.0₃0949'
$ws.Range("E29").Value = '  +10.43%  '

$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.02%  '

$ws.Range("D31").Value = "'522.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.62%  '

$ws.Range("E32").Value = '  +3.96%  '

$ws.Range("E33").Value = '  +2.43%  '

$ws.Range("D34").Value = "'5.64"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.52%  '

$ws.Range("E35").Value = '  +4.21%  '

$ws.Range("D36").Value = "'0.429"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.91%  '

$ws.Range("D37").Value = "'164.95"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.03%  '

$ws.Range("D38").Value = "'20.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +4.12%  '

$ws.Range("D39").Value = "'2.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.72%  '

$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("E41").Value = '  +0.08%  '

$ws.Range("D42").Value = "'42.18"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +6.50%  '

$ws.Range("D43").Value = "'165.62"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.66%  '

$ws.Range("E44").Value = '  +2.94%  '

$ws.Range("D45").Value = "'0.0621"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +6.45%  '

$ws.Range("D46").Value = "'23.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.67%  '

$ws.Range("E47").Value = '  +4.37%  '

$ws.Range("E48").Value = '  +3.60%  '

$ws.Range("D49").Value = "'0.0255"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.87%  '

$ws.Range("E50").Value = '  +2.76%  '

$ws.Range("D51").Value = "'19.54"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.98%  '
